$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.529.98"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").Value = "1.564.13"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("D4").Value = "'0.990"
$ws.Range("E4").Value = "  -1.66%  "

$ws.Range("D5").Value = "'210.79"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D8").Value = "'22.67"
$ws.Range("E8").Value = "  +2.67%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").Value = "1.787.31"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "1.563.54"
$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "27.508.65"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "'62.49"
$ws.Range("E17").Value = "  +0.95%  "

$ws.Range("D18").Value = "'225.55"
$ws.Range("E18").Value = "  +4.51%  "

$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").Value = "0.0₃0706"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("E21").Value = "  -1.68%  "

$ws.Range("D22").Value = "'4.12"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "'9.39"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").Value = "'1.95"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").Value = "'149.78"
$ws.Range("E25").Value = "  -2.62%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'15.19"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.108"
$ws.Range("E27").Value = "  +2.17%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.62"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").Value = "  -1.45%  "

$ws.Range("D30").Value = "'1.14"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").Value = "1.454.81"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("D34").Value = "'3.16"
$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("E35").Value = "  +2.65%  "

$ws.Range("E36").Value = "  +0.57%  "

$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").Value = "'0.0167"
$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("D40").Value = "'0.814"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").Value = "'5.72"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("E43").Value = "  -1.75%  "

$ws.Range("D44").Value = "'1.84"
$ws.Range("E44").Value = "  +5.63%  "

$ws.Range("D45").Value = "'0.972"
$ws.Range("E45").Value = "  -3.08%  "

$ws.Range("D46").Value = "'64.72"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "1.701.47"

$ws.Range("D48").Value = "'86.52"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").Value = "'0.0947"
$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("D51").Value = "'0.989"
$ws.Range("E51").Value = "  -1.71%  "
